$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 30.75
$ws.Range("I8").Value = 5.2
$ws.Range("K8").Value = 15.6
$ws.Range("M8").Value = 123.4

$ws.Range("H33").Value = 1005.7778
$ws.Range("I33").Value = 1005.7778
$ws.Range("K33").Value = 1005.7778
$ws.Range("M33").Value = -776.7778

$ws.Range("H98").Value = 715.44446
$ws.Range("I98").Value = 715.44446
$ws.Range("K98").Value = 715.44446
$ws.Range("M98").Value = 782.55554

$ws.Range("H122").Value = 715.44446
$ws.Range("I122").Value = 715.44446
$ws.Range("K122").Value = 2146.33338
$ws.Range("M122").Value = 303.66662

$ws.Range("H123").Value = 59999.77
$ws.Range("J123").Value = 59999.77
$ws.Range("L123").Value = 59999.77
$ws.Range("N123").Value = -69799.76999999999

$ws.Range("H138").Value = 2659.3635
$ws.Range("J138").Value = 3262.7048
$ws.Range("L138").Value = 9788.1144
$ws.Range("N138").Value = -20068.1144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11566.536
$ws.Range("I2").Value = 17178.445
$ws.Range("K2").Value = 17178.445
$ws.Range("M2").Value = -17065.445

$ws.Range("H32").Value = 11627.365
$ws.Range("I32").Value = 10747.807
$ws.Range("J32").Value = 19983.166
$ws.Range("K32").Value = 10747.807
$ws.Range("L32").Value = 19983.166
$ws.Range("M32").Value = -10460.807
$ws.Range("N32").Value = -20557.166

$ws.Range("H44").Value = 20000
$ws.Range("I44").Value = 20000
$ws.Range("K44").Value = 20000
$ws.Range("M44").Value = -19512

$ws.Range("H61").Value = 3489.7
$ws.Range("I61").Value = 2284.8572
$ws.Range("J61").Value = 4138.4614
$ws.Range("K61").Value = 2284.8572
$ws.Range("L61").Value = 4138.4614
$ws.Range("M61").Value = -2072.8572
$ws.Range("N61").Value = -4562.4614

$ws.Range("H74").Value = 1763.8
$ws.Range("I74").Value = 1763.8
$ws.Range("K74").Value = 1763.8
$ws.Range("M74").Value = -889.8

$ws.Range("H77").Value = 1763.8
$ws.Range("I77").Value = 1763.8
$ws.Range("K77").Value = 8819
$ws.Range("M77").Value = -4451

$ws.Range("H116").Value = 11566.536
$ws.Range("I116").Value = 17178.445
$ws.Range("K116").Value = 17178.445
$ws.Range("M116").Value = -14884.445

$ws.Range("H132").Value = 3196.1
$ws.Range("I132").Value = 2975.1538
$ws.Range("K132").Value = 8925.4614
$ws.Range("M132").Value = -6395.4614

$ws.Range("H136").Value = 3489.7
$ws.Range("I136").Value = 2284.8572
$ws.Range("J136").Value = 4138.4614
$ws.Range("K136").Value = 6854.571599999999
$ws.Range("L136").Value = 12415.3842
$ws.Range("M136").Value = -4304.571599999999
$ws.Range("N136").Value = -17515.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11566.536
$ws.Range("I3").Value = 17178.445
$ws.Range("K3").Value = 17178.445
$ws.Range("M3").Value = -17064.445

$ws.Range("H105").Value = 3443.25
$ws.Range("I105").Value = 2463.182
$ws.Range("K105").Value = 2463.182
$ws.Range("M105").Value = -716.1819999999998

$ws.Range("H134").Value = 2599.7837
$ws.Range("I134").Value = 1383.0667
$ws.Range("K134").Value = 4149.2001
$ws.Range("M134").Value = -1614.2001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3697.8823
$ws.Range("I31").Value = 1561.5
$ws.Range("K31").Value = 1561.5
$ws.Range("M31").Value = -1266.5

$ws.Range("H34").Value = 3697.8823
$ws.Range("I34").Value = 1561.5
$ws.Range("K34").Value = 1561.5
$ws.Range("M34").Value = -1359.5

$ws.Range("H99").Value = 8638023
$ws.Range("I99").Value = 3055992
$ws.Range("J99").Value = 11118926
$ws.Range("K99").Value = 3055992
$ws.Range("L99").Value = 11118926
$ws.Range("M99").Value = -3054494
$ws.Range("N99").Value = -11121922

$ws.Range("H126").Value = 8638023
$ws.Range("I126").Value = 3055992
$ws.Range("J126").Value = 11118926
$ws.Range("K126").Value = 9167976
$ws.Range("L126").Value = 33356778
$ws.Range("M126").Value = -9165506
$ws.Range("N126").Value = -33361718

$ws.Range("H141").Value = 153032.03
$ws.Range("J141").Value = 153032.03
$ws.Range("L141").Value = 153032.03
$ws.Range("N141").Value = -163392.03

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1112.7778
$ws.Range("I7").Value = 1260.2
$ws.Range("K7").Value = 3780.6
$ws.Range("M7").Value = -3668.6

$ws.Range("H45").Value = 4445
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 4445
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 13335
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -14399

$ws.Range("H64").Value = 2488.889
$ws.Range("J64").Value = 2488.889
$ws.Range("L64").Value = 7466.667
$ws.Range("N64").Value = -8006.667

$ws.Range("H67").Value = 2488.889
$ws.Range("J67").Value = 2488.889
$ws.Range("L67").Value = 7466.667
$ws.Range("N67").Value = -9338.667000000001

$ws.Range("H113").Value = 4221.1934
$ws.Range("I113").Value = 8150.4287
$ws.Range("J113").Value = 985.35297
$ws.Range("K113").Value = 24451.2861
$ws.Range("L113").Value = 2956.05891
$ws.Range("M113").Value = -22281.2861
$ws.Range("N113").Value = -7296.05891

$ws.Range("H131").Value = 3934.5908
$ws.Range("I131").Value = 3427.375
$ws.Range("K131").Value = 10282.125
$ws.Range("M131").Value = -5242.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 575.5238000000001
$ws.Range("I97").Value = 465.83334
$ws.Range("K97").Value = 465.83334
$ws.Range("M97").Value = 30.16665999999998

$ws.Range("H132").Value = 3016.0293
$ws.Range("I132").Value = 2814.5652
$ws.Range("J132").Value = 3437.2727
$ws.Range("K132").Value = 8443.695599999999
$ws.Range("L132").Value = 10311.8181
$ws.Range("M132").Value = -5913.695599999999
$ws.Range("N132").Value = -15371.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 3281.625
$ws.Range("J55").Value = 1120.6
$ws.Range("L55").Value = 1120.6
$ws.Range("N55").Value = -1466.6

$ws.Range("H82").Value = 2474.0908
$ws.Range("I82").Value = 1648.5217
$ws.Range("J82").Value = 4372.9
$ws.Range("K82").Value = 1648.5217
$ws.Range("L82").Value = 4372.9
$ws.Range("M82").Value = -1287.5217
$ws.Range("N82").Value = -5094.9

$ws.Range("H85").Value = 2474.0908
$ws.Range("I85").Value = 1648.5217
$ws.Range("J85").Value = 4372.9
$ws.Range("K85").Value = 1648.5217
$ws.Range("L85").Value = 4372.9
$ws.Range("M85").Value = -400.5217
$ws.Range("N85").Value = -6868.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3721.7144
$ws.Range("I62").Value = 3460.4
$ws.Range("K62").Value = 3460.4
$ws.Range("M62").Value = -2836.4

$ws.Range("H65").Value = 3721.7144
$ws.Range("I65").Value = 3460.4
$ws.Range("K65").Value = 17302
$ws.Range("M65").Value = -14182

$ws.Range("H94").Value = 60000
$ws.Range("J94").Value = 60000
$ws.Range("L94").Value = 60000
$ws.Range("N94").Value = -61802

$ws.Range("H100").Value = 664.1739
$ws.Range("I100").Value = 648.9091
$ws.Range("K100").Value = 1297.8182
$ws.Range("M100").Value = -756.8181999999999
